$d = $word.ActiveDocument
$table = $d.Tables(1)

# Helper-less direct calls (one per field) so each replacement is scoped to its
# own table cell; this keeps numerically-identical looking values in other
# cells/fields untouched, and MatchWholeWord keeps numbers like "27" from
# matching inside "127" etc.

# --- Record 1 (cell 1): "E31" / IKHSANUL QOMAR -> "E30" / I MADE SURYA D. ---
$cell1 = $table.Cell(1, 1)
$cell1.Range.Find.Execute("E31", $true, $true, $false, $false, $false, $true, 1, $false, "E30", 1)
$cell1.Range.Find.Execute("IKHSANUL QOMAR", $true, $false, $false, $false, $false, $true, 1, $false, "I MADE SURYA D.", 1)
$cell1.Range.Find.Execute("98", $true, $true, $false, $false, $false, $true, 1, $false, "102", 1)
$cell1.Range.Find.Execute("108", $true, $true, $false, $false, $false, $true, 1, $false, "127", 1)
$cell1.Range.Find.Execute("72", $true, $true, $false, $false, $false, $true, 1, $false, "76", 1)
$cell1.Range.Find.Execute("66", $true, $true, $false, $false, $false, $true, 1, $false, "39", 1)
$cell1.Range.Find.Execute("27", $true, $true, $false, $false, $false, $true, 1, $false, "29", 1)
$cell1.Range.Find.Execute("19", $true, $true, $false, $false, $false, $true, 1, $false, "23", 1)
$cell1.Range.Find.Execute("92", $true, $true, $false, $false, $false, $true, 1, $false, "97", 1)

# --- Record 2 (cell 2): "E68" / ARNOL KUSUMA -> "E38" / MUHAMMAD A. ---
$cell2 = $table.Cell(1, 2)
$cell2.Range.Find.Execute("E68", $true, $true, $false, $false, $false, $true, 1, $false, "E38", 1)
$cell2.Range.Find.Execute("ARNOL KUSUMA", $true, $false, $false, $false, $false, $true, 1, $false, "MUHAMMAD A.", 1)
$cell2.Range.Find.Execute("96", $true, $true, $false, $false, $false, $true, 1, $false, "85", 1)
$cell2.Range.Find.Execute("110", $true, $true, $false, $false, $false, $true, 1, $false, "106", 1)
$cell2.Range.Find.Execute("68", $true, $true, $false, $false, $false, $true, 1, $false, "65", 1)
$cell2.Range.Find.Execute("35", $true, $true, $false, $false, $false, $true, 1, $false, "34", 1)
$cell2.Range.Find.Execute("27", $true, $true, $false, $false, $false, $true, 1, $false, "26", 1)
$cell2.Range.Find.Execute("21", $true, $true, $false, $false, $false, $true, 1, $false, "19", 1)
$cell2.Range.Find.Execute("95", $true, $true, $false, $false, $false, $true, 1, $false, "94", 1)
$cell2.Range.Find.Execute("DP 2 TEKNIKA/15", $true, $false, $false, $false, $false, $true, 1, $false, "DP 2 NAUTIKA/15", 1)

# --- Record 3 (cell 3): "E73" / EKO PRASETYO -> "E98" / UTANG SUHAYA ---
$cell3 = $table.Cell(1, 3)
$cell3.Range.Find.Execute("E73", $true, $true, $false, $false, $false, $true, 1, $false, "E98", 1)
$cell3.Range.Find.Execute("EKO PRASETYO", $true, $false, $false, $false, $false, $true, 1, $false, "UTANG SUHAYA", 1)
$cell3.Range.Find.Execute("111", $true, $true, $false, $false, $false, $true, 1, $false, "94", 1)
$cell3.Range.Find.Execute("120", $true, $true, $false, $false, $false, $true, 1, $false, "102", 1)
$cell3.Range.Find.Execute("73", $true, $true, $false, $false, $false, $true, 1, $false, "72", 1)
$cell3.Range.Find.Execute("62", $true, $true, $false, $false, $false, $true, 1, $false, "33", 1)
# UC_5/UC_6 swap: UC_5 "20"->"26" then UC_6 "21"->"20" (sequenced so the
# second Find only matches the still-untouched UC_6 value).
$cell3.Range.Find.Execute("20", $true, $true, $false, $false, $false, $true, 1, $false, "26", 1)
$cell3.Range.Find.Execute("21", $true, $true, $false, $false, $false, $true, 1, $false, "20", 1)
$cell3.Range.Find.Execute("93", $true, $true, $false, $false, $false, $true, 1, $false, "91", 1)
